# The workbook's sheet data runs through row 129 (date serial 45685).
# Append two more daily rows (130, 131) that duplicate the values of the
# last existing row (129), only advancing the date in column A by one
# day for each new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 129 into rows 130 and 131, carrying over all formatting
# (e.g. the date style applied to column A).
$ws.Range("A129:J129").Copy($ws.Range("A130:J130"))
$ws.Range("A129:J129").Copy($ws.Range("A131:J131"))

# Set the correct dates for the two new rows (one day after each other).
$ws.Cells.Item(130, 1).Value2 = 45686
$ws.Cells.Item(131, 1).Value2 = 45687
